$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (shifts FightHeroPos columns etc. one to the right)
$ws.Columns("L:L").Insert()

# Populate the newly inserted "FightHeroCnfID" column (copy of FightHero metadata column)
$ws.Range("L1").Value = "FightHeroCnfID"
$ws.Range("L2").Value = "string"
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("L10").Value = "Hero"

# Update view state to match the post-edit selection
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("L3").Select()
